$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.370077013969421
$ws.Range("B1").Value = 1.384758234024048
$ws.Range("C1").Value = 1.488422155380249
$ws.Range("D1").Value = 2.132534027099609
$ws.Range("E1").Value = 4.273551464080811
